# Icris.Excel2Api.CoreWeb/sheets/test.xlsx - "Ironed out a few bugs, created 'visible' column."
$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("Input")
$wsProcess = $wb.Worksheets.Item("burp")
$wsOutput  = $wb.Worksheets.Item("Output")

# --- Rename 'burp' -> 'Process' ---
$wsProcess.Name = "Process"

# --- Input sheet: insert a new 'Visible' column (G) before the old 'Options' column ---
$wsInput.Columns.Item(7).Insert() | Out-Null

$wsInput.Range("G1").Value = "Visible"
$wsInput.Range("G2").Value = $true
$wsInput.Range("G3").Value = $true
$wsInput.Range("G4").Value = $true
$wsInput.Range("G5").Value = $true

# --- Fix a bug: Width of the object (D4) was 9, should be 8 ---
$wsInput.Range("D4").Value = 8

# --- Fix the VLOOKUP formula on the Output sheet to reference the renamed 'Process' sheet ---
$wsOutput.Range("C5").Formula = "=VLOOKUP(Input!D5,Process!A2:B5,2,FALSE)*C3/1000"

# --- Selections / active sheet, matching the final saved state ---
$wsProcess.Range("B6").Select() | Out-Null
$wsOutput.Range("C6").Select() | Out-Null
$wsInput.Activate() | Out-Null
$wsInput.Range("G6").Select() | Out-Null

$wb.Application.Calculate()
